# Correct some column translations
# The AT_2016/AT_2017/AT_2018 columns (E/F/G) had incorrectly been filled in
# with the same source-column names used by AT_2019 and later (columns
# H onward). Fix the rows where the AT_2016-2018 source schema actually
# differed: farm_id, crop_code, crop_name, organic, field_size.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: farm_id -> AT_2016/2017/2018 did not have this column, clear it
$ws.Range("E3:G3").ClearContents()

# Row 4: crop_code -> AT_2016/2017/2018 used "SNAR_CODE" instead of "snart_code"
$ws.Range("E4:G4").Value = "SNAR_CODE"

# Row 5: crop_name -> AT_2016/2017/2018 used "SNAR_BEZEICHNUNG" instead of "snart"
$ws.Range("E5:G5").Value = "SNAR_BEZEICHNUNG"

# Row 11: organic -> AT_2016/2017/2018 did not have this column, clear it
$ws.Range("E11:G11").ClearContents()

# Row 12: field_size -> AT_2016/2017/2018 used "SL_FLAECHE_BRUTTO_HA" (previously blank)
$ws.Range("E12:G12").Value = "SL_FLAECHE_BRUTTO_HA"

# Leave selection where the edit ended
$ws.Range("E16").Select()
